# Journal de travail - add new work-log entries (rows 28-32), matching the
# formatting (borders/number formats) already used by the existing rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the formatting of the last existing row (27) down across the five
# new rows first, so the new cells pick up the same styles (date format,
# borders, wrap text, etc.) as the rest of the table - including the blank,
# bordered F column.
$ws.Range("A27:F27").Copy() | Out-Null
$ws.Range("A28:F32").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 28 - 06.03.2023, semaine 5, Implémentation
$ws.Range("A28").Value = 44991
$ws.Range("B28").Value = 5
$ws.Range("C28").Value = "2,25"
$ws.Range("D28").Value = "Implémentation"
$ws.Range("E28").Value = "Afficher les valeurs de la base de données sur la page d'accueil"

# Row 29 - 06.03.2023, semaine 5, Implémentation
$ws.Range("A29").Value = 44991
$ws.Range("B29").Value = 5
$ws.Range("C29").Value = "2,25"
$ws.Range("D29").Value = "Implémentation"
$ws.Range("E29").Value = "Afficher les valeurs de la base de données sur la page d'accueil et modifier le front-end"

# Row 30 - 07.03.2023, semaine 6, Implémentation
$ws.Range("A30").Value = 44992
$ws.Range("B30").Value = 6
$ws.Range("C30").Value = "1,30"
$ws.Range("D30").Value = "Implémentation"
$ws.Range("E30").Value = "Créer un formulaire et afficher les nouvelles données"

# Row 31 - 09.03.2023, semaine 7, Implémentation
$ws.Range("A31").Value = 44994
$ws.Range("B31").Value = 7
$ws.Range("C31").Value = "1,30"
$ws.Range("D31").Value = "Implémentation"
$ws.Range("E31").Value = "Créer un formulaire et afficher les nouvelles données, modifier le fichier php de connection"

# Row 32 - 10.03.2023, semaine 8, Implémentation
$ws.Range("A32").Value = 44995
$ws.Range("B32").Value = 8
$ws.Range("C32").Value = "2,25"
$ws.Range("D32").Value = "Implémentation"
$ws.Range("E32").Value = "Finir le formulaire et mettre à jour les données sur swisscenter"

# Row heights: rows 28, 29, 31, 32 use the taller (wrapped) 30pt row like
# their siblings above; row 30 keeps the sheet's default height.
$ws.Range("A28:F28").RowHeight = 30
$ws.Range("A29:F29").RowHeight = 30
$ws.Range("A31:F31").RowHeight = 30
$ws.Range("A32:F32").RowHeight = 30

# Move the view down to the newly added rows and select the next empty row,
# matching the author's on-save cursor position.
$ws.Range("A34").Select() | Out-Null
